$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC01/TC05/TC06: refresh the sample CUID test customer number and let the
# dependent SQL-statement formula in H2 recompute against the new value.
$ws.Range("G2").Value = 10000700961

# Workbook is in manual calculation mode; force a recalc so the cached
# formula result (H2) is updated to match the new G2 input.
$excel.Calculate()
